$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.090784430503845
$ws.Range("B1").Value = 2.449971914291382
$ws.Range("C1").Value = 6.389114379882812
$ws.Range("D1").Value = 2.21128249168396
$ws.Range("E1").Value = 1.272898435592651
